# Update the 丽水-漫展信息 workbook to the next scraped snapshot.
# The event that used to be row 2 (丽水·CCAC动漫游戏嘉年华, 2024-07-20) has
# finished and drops off the list; every remaining event shifts up one row,
# a couple of "interested" counters (F) tick up slightly, and a brand new
# event (丽水·LZ栗子动漫游戏嘉年华, 2024-09-16) is appended at the end.
# Net effect: table shrinks from 9 data rows (A2:I10) to 7 data rows (A2:I8).

$wb = $excel.ActiveWorkbook

# New content for data rows 2-8 (row 1 is the header and is unchanged).
$rows = @(
    @{ A = 1; B = "2024-07-20"; C = "丽水·龙泉ACG动漫游戏博览会"; D = "南秦路1号望瓯·陶溪川直走200米左手边(7号楼) 望瓯陶溪川活动中心"; E = "2024.07.20 10:00-07.21 18:00"; F = 1718; G = 60;   H = "https://show.bilibili.com/platform/detail.html?id=86671"; I = "//i0.hdslb.com/bfs/openplatform/202406/LSorIT7S1717486817969.png" },
    @{ A = 2; B = "2024-07-27"; C = "丽水·thp01～风摄少微";          D = "大猷街 应星楼";                                                                   E = "2024.07.27 10:00-07.27 18:00"; F = 29;   G = 50;   H = "https://show.bilibili.com/platform/detail.html?id=87134"; I = "//i2.hdslb.com/bfs/openplatform/202406/JuvSmncN1717775885615.png" },
    @{ A = 3; B = "2024-07-27"; C = "丽水·第四届HP国风动漫游戏嘉年华";  D = "城北街798号 莱茵体育生活馆";                                                             E = "2024.07.27 08:30-07.27 17:00"; F = 479;  G = 65;   H = "https://show.bilibili.com/platform/detail.html?id=87305"; I = "//i2.hdslb.com/bfs/openplatform/202406/YUnPOKGV1718268952725.jpeg" },
    @{ A = 4; B = "2024-08-03"; C = "丽水·樱卡动漫游戏嘉年华";         D = "中东路848号(解放街交汇) 飞达国际大酒店";                                                     E = "2024.08.03 10:00-08.03 17:00"; F = 158;  G = 50;   H = "https://show.bilibili.com/platform/detail.html?id=87276"; I = "//i0.hdslb.com/bfs/openplatform/202406/bVp0Zg1B1718172430380.jpeg" },
    @{ A = 5; B = "2024-08-10"; C = "丽水·CCAC动漫七夕（回馈展）";     D = "中东路848号(解放街交汇) 飞达国际大酒店";                                                     E = "2024.08.10 09:00-08.10 17:00"; F = 79;   G = 29.9; H = "https://show.bilibili.com/platform/detail.html?id=86567"; I = "//i0.hdslb.com/bfs/openplatform/202405/tsOzbBRx1717015539538.png" },
    @{ A = 6; B = "2024-08-17"; C = "丽水·AEO纯白礼赞动漫嘉年华";      D = "城北街1001号 爱依·时尚婚宴中心";                                                         E = "2024.08.17 09:00-08.17 16:00"; F = 647;  G = 55;   H = "https://show.bilibili.com/platform/detail.html?id=86779"; I = "//i2.hdslb.com/bfs/openplatform/202406/MxJ3oNjt1717405405850.jpeg" },
    @{ A = 7; B = "2024-09-16"; C = "丽水·LZ栗子动漫游戏嘉年华";       D = "城北街798号 莱茵体育生活馆";                                                             E = "2024.09.16 09:30-09.16 17:00"; F = 410;  G = 65;   H = "https://show.bilibili.com/platform/detail.html?id=87480"; I = "//i1.hdslb.com/bfs/openplatform/202406/bATqcZhH1719285865931.jpeg" }
)

# Both the "展览" (exhibitions) sheet and the "全部类型" (all types) sheet
# carry this same table and both need the identical update.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the two oldest rows (old rows 9 and 10) first so the sheet's
    # used range shrinks to A1:I8 and the remaining rows (2-8) can simply
    # be overwritten in place with the refreshed data below.
    $ws.Range("A9:A10").EntireRow.Delete()

    for ($i = 0; $i -lt $rows.Length; $i++) {
        $r = $i + 2
        $row = $rows[$i]

        $ws.Cells.Item($r, 1).Value = $row.A

        # Column B holds plain "YYYY-MM-DD" text in the source data; writing
        # it straight through Value lets Excel reinterpret it as a date
        # serial, so force text mode for the write and then clear the
        # number-format override back off so the cell keeps the workbook's
        # default (unstyled) formatting, matching the original file.
        $bCell = $ws.Cells.Item($r, 2)
        $bCell.NumberFormat = "@"
        $bCell.Value = $row.B
        $bCell.ClearFormats()

        $ws.Cells.Item($r, 3).Value = $row.C
        $ws.Cells.Item($r, 4).Value = $row.D
        $ws.Cells.Item($r, 5).Value = $row.E
        $ws.Cells.Item($r, 6).Value = $row.F
        $ws.Cells.Item($r, 7).Value = $row.G
        $ws.Cells.Item($r, 8).Value = $row.H
        $ws.Cells.Item($r, 9).Value = $row.I
    }
}
